{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change summary (from the diff):\n//  1. In the \"Countryside\" formula paragraph, the edge tuple loses its\n//     \"poidsA\" field and the node tuple gains a trailing \"hungerS\" field:\n//       )(indiceA,fromS,toS,poidsA)(indiceS,imageS,posxS,posyS,valueS)\n//     becomes\n//       )(indiceA,fromS,toS)(indiceS,imageS,posxS,posyS,valueS,hungerS)\n//  2. In the data paragraph right below it, every edge's trailing\n//     \"poidsA\" value (always \"1\") is dropped, and every node (image)\n//     entry gets a new trailing \"hungerS\" number appended right after\n//     its existing valueS.\n//\n// We use Range.search()+insertText(..., \"Replace\") on the smallest\n// unique substrings so the surrounding content - including the\n// \"_GoBack\" bookmark that sits inside the data paragraph - is left\n// untouched.\n\nconst oldFormulaTail = \")(indiceA,fromS,toS,poidsA)(indiceS,imageS,posxS,posyS,valueS)\";\nconst newFormulaTail = \")(indiceA,fromS,toS)(indiceS,imageS,posxS,posyS,valueS,hungerS)\";\n\nconst oldData =\n  \"0 1 7 1 1 10 4 1 2 10 3 1 3 3 2 1 4 2 5 1 5 5 0 1 6 8 0 1 7 9 0 1 8 5 6 1 9 9 6 1 10 7 9 1 11 7 8 1 12 2 10 1 13 6 10 1 14 5 10 1 15 0 10 1 16 9 10 1 17 7 10 1 18 1 10 1 19 11 2 1 20 11 1 1 21 4 1 1 0serpent.jpg#454 360 20 1pin.jpg#780 410 50 2gland.jpg#120 400 50 3nitrate.jpg#100 650 80 4nitrate.jpg#850 650 80 5souris.jpg#330 450 30 6faucon.jpg#360 150 5 7scarabe.jpg#700 310 30 8salamandre.jpg#550 460 15 9moineau.jpg#560 170 20 10decomposer.jpg#454 620 80 11soleil.jpg#100 100 100\";\n\nconst newData =\n  \"0 1 7 1 10 4 2 10 3 3 3 2 4 2 5 5 5 0 6 8 0 7 9 0 8 5 6 9 9 6 10 7 9 11 7 8 12 2 10 13 6 10 14 5 10 15 0 10 16 9 10 17 7 10 18 1 10 19 11 2 20 11 1 21 4 1 0serpent.jpg#454 360 20 0.14 1pin.jpg#780 410 50 50 2gland.jpg#120 400 50 40 3nitrate.jpg#100 650 80 60 4nitrate.jpg#850 650 80 60 5souris.jpg#330 450 30 5 6faucon.jpg#360 150 5 1 7scarabe.jpg#700 310 30 0.5 8salamandre.jpg#550 460 15 2 9moineau.jpg#560 170 20 3 10decomposer.jpg#454 620 80 60 11soleil.jpg#100 100 100 0\";\n\nconst body = context.document.body;\n\n// --- 1. Patch the formula line -----------------------------------------\nconst formulaResults = body.search(oldFormulaTail, { matchCase: true });\nformulaResults.load(\"text\");\nawait context.sync();\n\nif (formulaResults.items.length === 0) {\n  throw new Error(\"Formula text to replace was not found: \" + oldFormulaTail);\n}\nformulaResults.items[0].insertText(newFormulaTail, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2. Patch the data line ----------------------------------------------\nconst dataResults = body.search(oldData, { matchCase: true });\ndataResults.load(\"text\");\nawait context.sync();\n\nif (dataResults.items.length === 0) {\n  throw new Error(\"Data text to replace was not found.\");\n}\ndataResults.items[0].insertText(newData, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change summary (from the diff):\n#  1. In the \"Countryside\" formula paragraph, the edge tuple loses its\n#     \"poidsA\" field and the node tuple gains a trailing \"hungerS\" field:\n#       )(indiceA,fromS,toS,poidsA)(indiceS,imageS,posxS,posyS,valueS)\n#     becomes\n#       )(indiceA,fromS,toS)(indiceS,imageS,posxS,posyS,valueS,hungerS)\n#  2. In the data paragraph right below it, every edge's trailing\n#     \"poidsA\" value (always \"1\") is dropped, and every node (image)\n#     entry gets a new trailing \"hungerS\" number appended right after\n#     its existing valueS.\n#\n# We use Find/Replace (wdReplaceOne) on the smallest unique substrings so\n# the surrounding content - including the \"_GoBack\" bookmark that sits\n# inside the data paragraph - is left untouched.\n\n$d = $word.ActiveDocument\n\n$oldFormulaTail = \")(indiceA,fromS,toS,poidsA)(indiceS,imageS,posxS,posyS,valueS)\"\n$newFormulaTail = \")(indiceA,fromS,toS)(indiceS,imageS,posxS,posyS,valueS,hungerS)\"\n\n$oldData = \"0 1 7 1 1 10 4 1 2 10 3 1 3 3 2 1 4 2 5 1 5 5 0 1 6 8 0 1 7 9 0 1 8 5 6 1 9 9 6 1 10 7 9 1 11 7 8 1 12 2 10 1 13 6 10 1 14 5 10 1 15 0 10 1 16 9 10 1 17 7 10 1 18 1 10 1 19 11 2 1 20 11 1 1 21 4 1 1 0serpent.jpg#454 360 20 1pin.jpg#780 410 50 2gland.jpg#120 400 50 3nitrate.jpg#100 650 80 4nitrate.jpg#850 650 80 5souris.jpg#330 450 30 6faucon.jpg#360 150 5 7scarabe.jpg#700 310 30 8salamandre.jpg#550 460 15 9moineau.jpg#560 170 20 10decomposer.jpg#454 620 80 11soleil.jpg#100 100 100\"\n\n$newData = \"0 1 7 1 10 4 2 10 3 3 3 2 4 2 5 5 5 0 6 8 0 7 9 0 8 5 6 9 9 6 10 7 9 11 7 8 12 2 10 13 6 10 14 5 10 15 0 10 16 9 10 17 7 10 18 1 10 19 11 2 20 11 1 21 4 1 0serpent.jpg#454 360 20 0.14 1pin.jpg#780 410 50 50 2gland.jpg#120 400 50 40 3nitrate.jpg#100 650 80 60 4nitrate.jpg#850 650 80 60 5souris.jpg#330 450 30 5 6faucon.jpg#360 150 5 1 7scarabe.jpg#700 310 30 0.5 8salamandre.jpg#550 460 15 2 9moineau.jpg#560 170 20 3 10decomposer.jpg#454 620 80 60 11soleil.jpg#100 100 100 0\"\n\n# --- 1. Patch the formula line -------------------------------------------\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Execute($oldFormulaTail, $false, $false, $false, $false, $false, $true, 1, $false, $newFormulaTail, 2)\n\n# --- 2. Patch the data line ------------------------------------------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute($oldData, $false, $false, $false, $false, $false, $true, 1, $false, $newData, 2)\n"}
